$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 51, shifting the existing row 51 (and below) down to 52.
$ws.Rows.Item(51).Insert()

# Row 52 now holds what used to be row 51 - it's already correct except for
# nothing (all original values moved down with the insert).

# Populate the new row 51 with the latest weekly price record.
$ws.Range("A51").Value = 10
$ws.Range("B51").Value = "Vega Modelo de Temuco"
$ws.Range("C51").Value = "La Araucanía"
$ws.Range("D51").Value = 44516
$ws.Range("E51").Value = 9
$ws.Range("F51").Value = 100112026
$ws.Range("G51").Value = "Haba"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 85
$ws.Range("K51").Value = 9000
$ws.Range("L51").Value = 9000
$ws.Range("M51").Value = 9000
$ws.Range("N51").Value = "$/saco 25 kilos"
$ws.Range("O51").Value = "Región del Maule"
$ws.Range("P51").Value = 360
$ws.Range("Q51").Value = 25
$ws.Range("R51").Value = "Hortaliza"

# Keep the date style (numFmt) on D51 consistent with the other date cells.
$ws.Range("D51").NumberFormat = $ws.Range("D50").NumberFormat
